# "reconsitution remove some unused file and code"
#
# The FilePath column (B) previously held a value only for row 2
# (pointing at an old/unused Ini\NFZoneServer\Scene\PioneerNoob\ folder).
# It is replaced with per-row references into Ini\Scene\<n>.xml for every
# data row (rows 2-7), and the new cells B3:B7 pick up a dedicated
# Text-formatted style (matching the existing "@" format already used by
# B2) with a small font tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (FilePath) values -----------------------------------------
$ws.Range("B2").Value = "../../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("B3").Value = "../../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("B4").Value = "../../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("B5").Value = "../../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("B6").Value = "../../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("B7").Value = "../../NFDataCfg/Ini/Scene/6.xml"

# --- Formatting for the newly-populated B3:B7 ---------------------------
# (B2 already carried the Text number format / style from before.)
$rng = $ws.Range("B3:B7")
$rng.NumberFormat = "@"
$rng.Font.Name = "宋体"
$rng.Font.Size = 11
$rng.Font.Family = 3

# --- Selection / scroll position -----------------------------------------
# Previously the view was scrolled to show column C first with E2
# selected; now it shows column A again with B5 selected.
$ws.Range("B5").Select() | Out-Null
